$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows per the diff
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -8
$ws.Range("F6").Value = -5
$ws.Range("F10").Value = -6
$ws.Range("F11").Value = -2
$ws.Range("F16").Value = -6
$ws.Range("F18").Value = 4
$ws.Range("F21").Value = -2
$ws.Range("F24").Value = -4
$ws.Range("F26").Value = -9
$ws.Range("F29").Value = 2
